{"js": "// Add factorial operation documentation to calculator_documentation.docx\n//\n// 1) Table of Contents entry \"4. Feature Specifications\" gets a line break\n//    followed by a new bullet line describing the Factorial Operation.\n// 2) The \"Menu Structure\" table (Option / Function) gets a new row:\n//    Factorial | factorial\n\n// --- 1. Amend the \"4. Feature Specifications\" Table-of-Contents entry ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst tocTarget = \"4. Feature Specifications\";\nlet tocParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === tocTarget) {\n    tocParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!tocParagraph) {\n  throw new Error('Could not find the \"4. Feature Specifications\" paragraph.');\n}\n\nconst bulletText =\n  \"\\u000b\\u2022 Factorial Operation: Calculates the factorial of an integer using \" +\n  \"recursion or iteration and logs each step in history if applicable. Includes \" +\n  \"error handling for non-integer inputs, with appropriate feedback to users.\";\n\n// Insert a manual line break (\\u000b) followed by the new bullet text at the\n// end of the existing run/paragraph, keeping it all inside the same run\n// (mirrors <w:t>...</w:t><w:br/><w:t>...</w:t> within one <w:r>).\ntocParagraph.insertText(bulletText, Word.InsertLocation.end);\n\n// --- 2. Append a \"Factorial\" row to the Menu Structure table ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst menuTable = tables.items[0];\nmenuTable.addRows(Word.InsertLocation.end, 1, [[\"Factorial\", \"factorial\"]]);\n\nawait context.sync();\n", "ps1": "# Add factorial operation documentation to calculator_documentation.docx\n#\n# 1) Table of Contents entry \"4. Feature Specifications\" gets a line break\n#    followed by a new bullet line describing the Factorial Operation.\n# 2) The \"Menu Structure\" table (Option / Function) gets a new row:\n#    Factorial | factorial\n\n$d = $word.ActiveDocument\n\n# --- 1. Amend the \"4. Feature Specifications\" Table-of-Contents entry ---\n$bulletText = [char]0x2022 + \" Factorial Operation: Calculates the factorial of an integer using recursion or iteration and logs each step in history if applicable. Includes error handling for non-integer inputs, with appropriate feedback to users.\"\n\nforeach ($p in $d.Paragraphs) {\n  $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($txt -eq \"4. Feature Specifications\" -and $p.Style.NameLocal -eq \"List Number\") {\n    $p.Range.InsertAfter([char]11 + $bulletText)\n    break\n  }\n}\n\n# --- 2. Append a \"Factorial\" row to the Menu Structure table ---\n$tbl = $d.Tables.Item(1)\n$newRow = $tbl.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"Factorial\"\n$newRow.Cells.Item(2).Range.Text = \"factorial\"\n"}
